# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51) with
# refreshed values scraped on Sun Nov  5 17:34:57 UTC 2023.
#
# A handful of coins changed rank position between scrapes, so for rows
# 16/17, 25/26 and 41/42 the Coin name + Link + Price + Volume are all
# rewritten together (not just a single cell).
#
# Several Price (column D) values are plain decimal numbers (e.g. "244.47").
# Setting .Value on such a string lets Excel's COM layer auto-convert it to
# a numeric cell, which would silently flip the cell from a text cell
# (t="inlineStr"/t="s") to a numeric cell (t="n") - not what the source
# data looks like. To keep those cells as text (matching the original
# file), we force NumberFormat to Text ("@") before writing the value and
# then restore the style to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.405.94"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.908.47"
$ws.Range("E3").Value = "  +2.77%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("E6").Value = "  +6.11%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.67"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  +6.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0717"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "2.184.25"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.19%  "
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.888.26"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "35.396.78"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("D20").Value = "0.0₃0823"
$ws.Range("E20").Value = "  +3.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "240.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +24.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.71%  "
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.935"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.17%  "
$ws.Range("E35").Value = "  +3.14%  "
$ws.Range("E36").Value = "  -3.82%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("E40").Value = "  +4.45%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.00%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0652"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "90.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "1.344.75"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +39.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "2.094.36"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0698"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.75%  "
